# ConfigurationFiles/2-Score_weight_Cer.xlsx
#
# Commit: "Start adding the filtes for the Ceramides where the MS/MS are
# checked if are [M+H]+ or [M+H-H2O]+"
#
# Functional change captured by this script:
#   - Rename the single worksheet from the generic "Sheet 1" to "Cer"
#     (the sheet now specifically holds the Ceramides scoring weights).
#   - Leave the cursor/selection on C25, matching where the author was
#     last working on the sheet when they saved it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook only has one worksheet ("Sheet 1") -> rename it to "Cer".
$ws.Name = "Cer"

# Move the active selection to C25 (was A40 before the edit).
$ws.Range("C25").Select()
